# Refresh cryptos list values (price/volume snapshot) to match the
# latest scrape. Price cells that look numeric are written with a
# leading apostrophe so Excel stores them as text (matching the
# original inline-string cells, e.g. keeping trailing zeros like
# "0.550" or "5.00" instead of collapsing them to 0.55 / 5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.077.72"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.916.16"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'202.64"
$ws.Range("E5").Value = "  +8.27%  "
$ws.Range("D6").Value = "'597.99"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D8").Value = "'0.550"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.197"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "2.915.08"
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("D11").Value = "'0.436"
$ws.Range("E11").Value = "  +17.75%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "'4.89"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "3.450.54"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "75.895.12"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "'27.90"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").Value = "'0.0000190"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "2.912.02"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'12.93"
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'8.85"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'372.44"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'2.33"
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("D23").Value = "'4.31"
$ws.Range("E23").Value = "  +6.30%  "
$ws.Range("D24").Value = "'71.38"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "3.053.78"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "'4.23"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "'9.73"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "'0.0000108"
$ws.Range("E29").Value = "  +3.79%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'502.72"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'7.72"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "'164.24"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").Value = "'20.18"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").Value = "'0.105"
$ws.Range("E39").Value = "  +23.59%  "
$ws.Range("D40").Value = "'0.113"
$ws.Range("E40").Value = "  -4.18%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'181.63"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'0.360"
$ws.Range("E43").Value = "  +6.64%  "
$ws.Range("D44").Value = "'5.00"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "'1.65"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'40.03"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'2.35"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'0.573"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'3.71"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'22.27"
$ws.Range("E51").Value = "  +6.60%  "
